# chore: update Sheets via scheduled runner
# Refresh cached profit-calculation figures (columns H-N) across the
# per-job worksheets to reflect the latest market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 458.875
$ws.Range("I9").Value = 650.2
$ws.Range("J9").Value = 140
$ws.Range("K9").Value = 650.2
$ws.Range("L9").Value = 140
$ws.Range("M9").Value = -481.2
$ws.Range("N9").Value = -478
$ws.Range("H28").Value = 7690.35
$ws.Range("I28").Value = 221.82353
$ws.Range("J28").Value = 50012
$ws.Range("K28").Value = 221.82353
$ws.Range("L28").Value = 50012
$ws.Range("M28").Value = 263.17647
$ws.Range("N28").Value = -50982
$ws.Range("H32").Value = 972.1875
$ws.Range("I32").Value = 1070.25
$ws.Range("J32").Value = 939.5
$ws.Range("K32").Value = 1070.25
$ws.Range("L32").Value = 939.5
$ws.Range("M32").Value = -744.25
$ws.Range("N32").Value = -1591.5
$ws.Range("H40").Value = 1832.5
$ws.Range("J40").Value = 1955.5
$ws.Range("L40").Value = 1955.5
$ws.Range("N40").Value = -2305.5
$ws.Range("H41").Value = 159.28572
$ws.Range("I41").Value = 250
$ws.Range("J41").Value = 144.16667
$ws.Range("K41").Value = 250
$ws.Range("L41").Value = 144.16667
$ws.Range("M41").Value = 190
$ws.Range("N41").Value = -1024.16667
$ws.Range("H43").Value = 1841.4286
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1841.4286
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1841.4286
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1979.4286
$ws.Range("H53").Value = 239.53847
$ws.Range("J53").Value = 363
$ws.Range("L53").Value = 363
$ws.Range("N53").Value = -1637
$ws.Range("H62").Value = 2983
$ws.Range("J62").Value = 3202.2
$ws.Range("L62").Value = 3202.2
$ws.Range("N62").Value = -4450.2
$ws.Range("H65").Value = 2983
$ws.Range("J65").Value = 3202.2
$ws.Range("L65").Value = 16011
$ws.Range("N65").Value = -22251
$ws.Range("H98").Value = 999.3200000000001
$ws.Range("I98").Value = 554.6111
$ws.Range("J98").Value = 2142.8572
$ws.Range("K98").Value = 554.6111
$ws.Range("L98").Value = 2142.8572
$ws.Range("M98").Value = 943.3889
$ws.Range("N98").Value = -5138.8572
$ws.Range("H113").Value = 2774
$ws.Range("I113").Value = 2700
$ws.Range("K113").Value = 2700
$ws.Range("M113").Value = 554
$ws.Range("H116").Value = 15459.25
$ws.Range("I116").Value = 27101.25
$ws.Range("J116").Value = 3817.25
$ws.Range("K116").Value = 27101.25
$ws.Range("L116").Value = 3817.25
$ws.Range("M116").Value = -23659.25
$ws.Range("N116").Value = -10701.25
$ws.Range("H122").Value = 999.3200000000001
$ws.Range("I122").Value = 554.6111
$ws.Range("J122").Value = 2142.8572
$ws.Range("K122").Value = 1663.8333
$ws.Range("L122").Value = 6428.571599999999
$ws.Range("M122").Value = 786.1667000000002
$ws.Range("N122").Value = -11328.5716
$ws.Range("H132").Value = 4903.2812
$ws.Range("I132").Value = 4600.2593
$ws.Range("K132").Value = 13800.7779
$ws.Range("M132").Value = -11270.7779
$ws.Range("H138").Value = 2351.966
$ws.Range("I138").Value = 2074.1365
$ws.Range("K138").Value = 6222.4095
$ws.Range("M138").Value = -1082.4095
$ws.Range("H141").Value = 4737.2
$ws.Range("I141").Value = 2109.8262
$ws.Range("J141").Value = 13370
$ws.Range("K141").Value = 6329.4786
$ws.Range("L141").Value = 40110
$ws.Range("M141").Value = -1149.4786
$ws.Range("N141").Value = -50470

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 778
$ws.Range("I21").Value = 778
$ws.Range("K21").Value = 778
$ws.Range("M21").Value = -404
$ws.Range("H32").Value = 409558.97
$ws.Range("I32").Value = 446172.66
$ws.Range("K32").Value = 446172.66
$ws.Range("M32").Value = -445885.66
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2126
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -10632
$ws.Range("N77").ClearContents()
$ws.Range("H97").Value = 863.375
$ws.Range("I97").Value = 840.93335
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 840.93335
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -344.93335
$ws.Range("N97").Value = -2192
$ws.Range("H122").Value = 54406.26
$ws.Range("I122").Value = 72988.57000000001
$ws.Range("J122").Value = 2375.8
$ws.Range("K122").Value = 218965.71
$ws.Range("L122").Value = 7127.400000000001
$ws.Range("M122").Value = -216515.71
$ws.Range("N122").Value = -12027.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1226.7858
$ws.Range("I94").Value = 852.8889
$ws.Range("J94").Value = 1899.8
$ws.Range("K94").Value = 852.8889
$ws.Range("L94").Value = 1899.8
$ws.Range("M94").Value = -401.8889
$ws.Range("N94").Value = -2801.8
$ws.Range("H134").Value = 2835.3667
$ws.Range("I134").Value = 2644
$ws.Range("J134").Value = 3085.6155
$ws.Range("K134").Value = 7932
$ws.Range("L134").Value = 9256.8465
$ws.Range("M134").Value = -5397
$ws.Range("N134").Value = -14326.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6954.65
$ws.Range("I31").Value = 1613.6666
$ws.Range("J31").Value = 9243.643
$ws.Range("K31").Value = 1613.6666
$ws.Range("L31").Value = 9243.643
$ws.Range("M31").Value = -1318.6666
$ws.Range("N31").Value = -9833.643
$ws.Range("H34").Value = 6954.65
$ws.Range("I34").Value = 1613.6666
$ws.Range("J34").Value = 9243.643
$ws.Range("K34").Value = 1613.6666
$ws.Range("L34").Value = 9243.643
$ws.Range("M34").Value = -1411.6666
$ws.Range("N34").Value = -9647.643
$ws.Range("H99").Value = 1619.25
$ws.Range("I99").Value = 723.1429000000001
$ws.Range("K99").Value = 723.1429000000001
$ws.Range("M99").Value = 774.8570999999999
$ws.Range("H126").Value = 1619.25
$ws.Range("I126").Value = 723.1429000000001
$ws.Range("K126").Value = 2169.4287
$ws.Range("M126").Value = 300.5712999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1270.2727
$ws.Range("I70").Value = 997.3
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 2991.9
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -2676.9
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 1270.2727
$ws.Range("I73").Value = 997.3
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 2991.9
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -1899.9
$ws.Range("N73").Value = -14184
$ws.Range("H113").Value = 852.9429
$ws.Range("J113").Value = 1108.238
$ws.Range("L113").Value = 3324.714
$ws.Range("N113").Value = -7664.714
$ws.Range("H114").Value = 1276.3125
$ws.Range("J114").Value = 2685.4285
$ws.Range("L114").Value = 8056.2855
$ws.Range("N114").Value = -14564.2855
$ws.Range("H117").Value = 1459.6
$ws.Range("J117").Value = 3225
$ws.Range("L117").Value = 9675
$ws.Range("N117").Value = -16559
$ws.Range("H129").Value = 1817.7
$ws.Range("I129").Value = 530
$ws.Range("J129").Value = 2044.9412
$ws.Range("K129").Value = 1590
$ws.Range("L129").Value = 6134.8236
$ws.Range("M129").Value = 3410
$ws.Range("N129").Value = -16134.8236
$ws.Range("H131").Value = 1009.43243
$ws.Range("J131").Value = 1079.3438
$ws.Range("L131").Value = 3238.0314
$ws.Range("N131").Value = -13318.0314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3055.111
$ws.Range("I132").Value = 2015.2307
$ws.Range("K132").Value = 6045.6921
$ws.Range("M132").Value = -3515.6921
$ws.Range("H136").Value = 8773536
$ws.Range("I136").Value = 1315.6923
$ws.Range("K136").Value = 3947.0769
$ws.Range("M136").Value = -1397.0769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2817.697
$ws.Range("I136").Value = 2858.3684
$ws.Range("J136").Value = 2762.5
$ws.Range("K136").Value = 8575.1052
$ws.Range("L136").Value = 8287.5
$ws.Range("M136").Value = -6025.1052
$ws.Range("N136").Value = -13387.5
